$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Long multi-paragraph article bodies (Alt+Enter style line breaks within the cell)
$hondaContent = @"
燃料ポンプの不具合で走行中にエンストする恐れがあるとして、ホンダは113万台余りのリコールを国土交通省に届け出ました。
リコールの対象は2017年6月から2020年9月までに製造されたホンダの「N−BOX」や「フィット」「FREED」など113万8000台余りです。
国交省によりますと、内部が変形して燃料ポンプが動かなくなり、最悪の場合は走行中にエンストする恐れがあるとしています。
これまで422件の不具合が確認されていますが、事故は起きていないということです。ホンダが同様の不具合でリコールを届け出るのは6回目で、アメリカやイギリスで生産された車を含めると合わせて約171万台に上ります。
同様の事象は他社でも確認されていて、一連の不具合によるリコールの届け出は2020年以降、合わせて380万台を超えています。
"@

$crawfishContent = @"
人口の増加や気候変動などによる世界的な食料不足への解決策として、東京農業大学がザリガニを新たな食資源として活用する研究を発表しました。
東京農業大学は、ケニアの大学や企業と共にアメリカザリガニの食用資源化に向けた研究を行っています。
   日本ではアメリカザリガニは生態系を脅かすとして今年6月から条件付特定外来生物に指定されていて、販売を目的とした飼育が禁止されています。
しかし、世界では養殖されているザリガニの99％以上がアメリカザリガニで、主に中国やアメリカで食用にされています。
高たんぱくで成長スピードが速く、生産コストも低いことから栄養不足が深刻なアフリカでも、食文化に取り入れるための研究を続けていくとしています。
東京農業大学の江口文陽学長は、「砂抜きや臭い抜きをすればザリガニも非常に素晴らしい食資源になると感じた。世界各国に食用生産、食用の活用がつながっていくことに大きな意味がある」と話しています。
"@

$putinContent = @"
ロシアのプーチン大統領は８日、来年３月１７日に実施予定の大統領選に出馬する意向を表明しました。通算５選目を目指す出馬で、当選すれば少なくとも２０３０年まで権力の座を維持することになります。
選挙は併合したウクライナのドネツク、ルハンスク、ヘルソン、ザポリージャの各州で行われる初めての大統領選ともなります。ロシアの中央選管によれば、これら４州では３月１５～１７日に「戸別訪問による投票」を実施するということです。
以前、親ロシア派の当局者がこれらの州で地方選を行った際には、国際社会から茶番と非難する声が上がっていました。
"@

$muskContent = @"
Ｘ（旧ツイッター）を所有する米起業家イーロン・マスク氏は７日、Ｘから広告を引き揚げた米ディズニーのボブ・アイガー最高経営責任者（ＣＥＯ）の決定を批判し、アイガー氏は強制解任されるべきだと明言した。
マスク氏はアイガー氏についてＸに「直ちに解任されるべきだ」と書き込み、「（創業者の）ウォルト・ディズニーはボブ（・アイガー氏）が自分の会社にやったことを見て、墓の中で安らかに眠れないだろう」と述べた。
ディズニーの代理人にコメントを求めたものの、現時点で応じていない。ただ、アイガー氏は最高経営責任者（ＣＥＯ）としての最初の任期中、「スター・ウォーズ」シリーズやマーベル・スタジオ、ピクサーの巧みな買収を通じディズニーを巨大エンタメ企業に押し上げた手腕が広く評価されている。
ディズニーは先月、他の多くの大企業と同様にＸへの広告出稿を停止した。マスク氏が白人至上主義者の間で受けの良い反ユダヤ主義陰謀論を支持したためだ。企業が相次ぎＸとの関係を断ったことを受け、マスク氏は先週、問題の投稿について暗に謝罪したものの、それと同時に広告購入を控える企業を汚い言葉を交えて罵倒した。
米紙ニューヨーク・タイムズのイベントに登壇したマスク氏は、「クソくらえ」という言葉を連発した。
アイガー氏はこれに先立ち、ディズニーはマスク氏を理由にＸへの広告出稿停止を決めたと説明しており、イベントではマスク氏がアイガー氏を名指しして罵詈（ばり）雑言を浴びせる場面もあった。
マスク氏は２０２２年後半にＸを買収して以降、ヘイトスピーチ（憎悪表現）や誤情報、陰謀論の増加を招く決定を相次ぎ下している。
"@

$putinTitle = "プーチン氏、２４年大統領選への出馬を表明"
$muskTitle = "マスク氏、ディズニーＣＥＯの「解任」要求　Ｘからの広告引き揚げを批判"

# --- Row 2: crawfish (food-shortage) article -- fix the content cell to use
# real line breaks instead of the ideographic-space run-on paragraph
$ws.Range("B2").Value = $crawfishContent

# --- Row 3: Honda recall article -- same line-break fix
$ws.Range("B3").Value = $hondaContent

# --- Row 4: new article -- Putin 2024 election bid ---
$ws.Range("A4").Value = $putinTitle
$ws.Range("B4").Value = $putinContent
$ws.Range("C4").Value = "Hv3Y6DjOO6J6bRNawzoXEvRKA31292023.jpeg"
$ws.Range("D4").Value = "329486f68b2b4177e0301e5c7256e277.mp3"

# --- Row 5: new article -- Musk vs Disney CEO ---
$ws.Range("A5").Value = $muskTitle
$ws.Range("B5").Value = $muskContent
$ws.Range("C5").Value = "affjttKJWeoxUTE8OhAHDKBJdk1292023.jpg"
$ws.Range("D5").Value = "329486f68b2b4177e0301e5c7266e600.mp3"

# Row heights (Excel recalculates these from the wrapped paragraph length;
# pin them to the values produced when this workbook was authored)
$ws.Rows(2).RowHeight = 409.5
$ws.Rows(3).RowHeight = 342
$ws.Rows(4).RowHeight = 252
$ws.Rows(5).RowHeight = 409.5

# Clear out the stray #NAME? formula that had been left in D20
$ws.Range("D20").ClearContents()

# Leave the selection where the author left it while reviewing the new rows
$ws.Range("F5").Select()
